$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers need an explicit Text
# number format first, otherwise Excel auto-converts them to numeric values
# (e.g. "2.00" -> 2), which would not match the source data (text strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "51.780.37"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.777.36"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "357.89"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "109.54"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("D7").Value = "0.563"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "39.94"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").Value = "7.58"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "3.215.35"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "2.754.01"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "0.932"
$ws.Range("E17").Value = "  +4.09%  "
$ws.Range("D18").Value = "51.745.73"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "7.41"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "3.11"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").Value = "12.98"
$ws.Range("E21").Value = "  -3.82%  "
$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").Value = "274.25"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "70.16"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").Value = "26.63"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").Value = "0.144"
$ws.Range("E29").Value = "  +3.39%  "
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "0.0463"
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").Value = "51.39"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "33.85"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "5.68"
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("D35").Value = "0.0842"
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("D36").Value = "5.23"
$ws.Range("E36").Value = "  +7.20%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "18.05"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "2.00"
$ws.Range("E40").Value = "  -3.92%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "2.53"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.115"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("D44").Value = "121.87"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("D45").Value = "22.07"
$ws.Range("E45").Value = "  -7.20%  "
$ws.Range("D46").Value = "2.066.43"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "3.23"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("E48").Value = "  -5.90%  "
$ws.Range("D49").Value = "5.68"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "0.929"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("E51").Value = "  +0.04%  "
